$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 266, pushing all existing rows (and the
# ones below) down by two. This mirrors the diff: dimension grows from
# A1:R348 to A1:R350, and the data previously on rows 266+ now lives on
# rows 268+.
$ws.Rows("266:267").Insert()

# --- New row 266 ---
$ws.Range("A266").Value = 3
$ws.Range("B266").Value = "Femacal de La Calera"
$ws.Range("C266").Value = "Coquimbo"
$ws.Range("D266").Value = 44588
$ws.Range("E266").Value = 5
$ws.Range("F266").Value = 100112032
$ws.Range("G266").Value = "Zapallo italiano"
$ws.Range("H266").Value = "Sin especificar"
$ws.Range("I266").Value = "Primera"
$ws.Range("J266").Value = 110
$ws.Range("K266").Value = 4000
$ws.Range("L266").Value = 4500
$ws.Range("M266").Value = 4273
$ws.Range("N266").Value = "`$/caja 36 unidades"
$ws.Range("O266").Value = "Provincia de Quillota"
$ws.Range("P266").Value = 119
$ws.Range("Q266").Value = 36
$ws.Range("R266").Value = "Hortaliza"

# --- New row 267 ---
$ws.Range("A267").Value = 3
$ws.Range("B267").Value = "Femacal de La Calera"
$ws.Range("C267").Value = "Coquimbo"
$ws.Range("D267").Value = 44588
$ws.Range("E267").Value = 5
$ws.Range("F267").Value = 100112032
$ws.Range("G267").Value = "Zapallo italiano"
$ws.Range("H267").Value = "Sin especificar"
$ws.Range("I267").Value = "Primera"
$ws.Range("J267").Value = 220
$ws.Range("K267").Value = 9000
$ws.Range("L267").Value = 9500
$ws.Range("M267").Value = 9277
$ws.Range("N267").Value = "`$/caja 70 unidades"
$ws.Range("O267").Value = "Provincia de Quillota"
$ws.Range("P267").Value = 133
$ws.Range("Q267").Value = 70
$ws.Range("R267").Value = "Hortaliza"
